{"js": "const replacements = [\n  [\"72\u00d729=\", \"72\u00d753=\"],\n  [\"93\u00d758=\", \"55\u00d731=\"],\n  [\"79\u00d713=\", \"86\u00d789=\"],\n  [\"76\u00d772=\", \"58\u00d715=\"],\n  [\"80\u00d750=\", \"96\u00d723=\"],\n  [\"43\u00d769=\", \"30\u00d777=\"],\n  [\"58\u00d772=\", \"19\u00d713=\"],\n  [\"47\u00d757=\", \"48\u00d764=\"],\n  [\"36\u00d783=\", \"70\u00d749=\"],\n  [\"71\u00d743=\", \"73\u00d790=\"],\n  [\"45\u00d711=\", \"67\u00d748=\"],\n  [\"69\u00d772=\", \"55\u00d738=\"],\n  [\"96\u00d719=\", \"57\u00d761=\"],\n  [\"91\u00d717=\", \"52\u00d711=\"],\n  [\"88\u00d726=\", \"37\u00d774=\"],\n  [\"34\u00d758=\", \"24\u00d777=\"],\n  [\"97\u00d713=\", \"78\u00d769=\"],\n  [\"99\u00d789=\", \"35\u00d723=\"],\n  [\"98\u00d793=\", \"75\u00d749=\"],\n  [\"99\u00d777=\", \"83\u00d791=\"],\n  [\"48\u00d759=\", \"18\u00d787=\"],\n  [\"19\u00d784=\", \"16\u00d789=\"],\n  [\"42\u00d797=\", \"92\u00d741=\"],\n  [\"55\u00d791=\", \"49\u00d779=\"],\n  [\"59\u00d729=\", \"72\u00d760=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"72\u00d729=\", \"72\u00d753=\"),\n    @(\"93\u00d758=\", \"55\u00d731=\"),\n    @(\"79\u00d713=\", \"86\u00d789=\"),\n    @(\"76\u00d772=\", \"58\u00d715=\"),\n    @(\"80\u00d750=\", \"96\u00d723=\"),\n    @(\"43\u00d769=\", \"30\u00d777=\"),\n    @(\"58\u00d772=\", \"19\u00d713=\"),\n    @(\"47\u00d757=\", \"48\u00d764=\"),\n    @(\"36\u00d783=\", \"70\u00d749=\"),\n    @(\"71\u00d743=\", \"73\u00d790=\"),\n    @(\"45\u00d711=\", \"67\u00d748=\"),\n    @(\"69\u00d772=\", \"55\u00d738=\"),\n    @(\"96\u00d719=\", \"57\u00d761=\"),\n    @(\"91\u00d717=\", \"52\u00d711=\"),\n    @(\"88\u00d726=\", \"37\u00d774=\"),\n    @(\"34\u00d758=\", \"24\u00d777=\"),\n    @(\"97\u00d713=\", \"78\u00d769=\"),\n    @(\"99\u00d789=\", \"35\u00d723=\"),\n    @(\"98\u00d793=\", \"75\u00d749=\"),\n    @(\"99\u00d777=\", \"83\u00d791=\"),\n    @(\"48\u00d759=\", \"18\u00d787=\"),\n    @(\"19\u00d784=\", \"16\u00d789=\"),\n    @(\"42\u00d797=\", \"92\u00d741=\"),\n    @(\"55\u00d791=\", \"49\u00d779=\"),\n    @(\"59\u00d729=\", \"72\u00d760=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Execute(\n        $old,       # FindText\n        $false,     # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $new,       # ReplaceWith\n        2           # Replace (wdReplaceAll)\n    )\n}\n"}
